$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.994.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.96%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.263.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.40%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.77%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.42%  "

# Row 7
$ws.Range("E7").Value = "  -1.26%  "

# Row 8
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.555"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.40%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.19%  "

# Row 11
$ws.Range("E11").Value = "  +0.60%  "

# Row 12
$ws.Range("E12").Value = "  -0.31%  "

# Row 13
$ws.Range("E13").Value = "  -2.05%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.608.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.18%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.863"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.44%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.19%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.269.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.90%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.898.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.84%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -9.88%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0985"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.63%  "

# Row 21
$ws.Range("E21").Value = "  +0.30%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.71"
$ws.Range("D22").Style = "Normal"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.81%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.30%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.07%  "

# Row 26
$ws.Range("E26").Value = "  +0.15%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.91%  "

# Row 28
$ws.Range("E28").Value = "  -3.54%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.73%  "

# Row 30
$ws.Range("E30").Value = "  -1.95%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.62%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.54%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0854"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.91%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.21%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.116"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.14%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.15%  "

# Row 37
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.71%  "

# Row 38
$ws.Range("E38").Value = "  -2.42%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +19.46%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.64%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.63%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0316"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.94%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.795.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.14%  "

# Row 45
$ws.Range("B45").Value = "ordi"
$ws.Range("C45").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "75.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.42%  "

# Row 46
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.199"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.29%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "82.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.04%  "

# Row 48
$ws.Range("E48").Value = "  -2.14%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.29%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "58.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.65%  "

# Row 51
$ws.Range("E51").Value = "  +5.39%  "
